# Apply the data changes described by the commit:
#   "Add ParaBank base, pages, tests and updated Excel data"
#
# Concretely, on the ParaBank_RegistrationForm sheet the scenario result
# changes from a "PASS" (account created, message shown on one line wrapped
# across 3 physical lines) to a "FAIL" (message normalized to a single
# line, no longer requiring wrap formatting).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ParaBank_RegistrationForm")

# Actual Results (M2): replace the old, hard-wrapped message with the
# single-line version of the same text, and drop the wrap-text styling
# that is no longer needed since the text no longer contains line breaks.
$ws.Range("M2").Value = "Your account was created successfully. You are now logged in."
$ws.Range("M2").Style = "Normal"

# Status (N2): the scenario now reports a failure instead of a pass.
$ws.Range("N2").Value = "FAIL"
